$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add notes (typed in this order so shared-string indices match: 35, 36, 37)
$ws.Range("G35").Value = "Added in a no checkpoint mode"
$ws.Range("G37").Value = "Added a death counter"
$ws.Range("G33").Value = "Just to add the points of the death counter being added"

# Tick the checkboxes (linked cells) for rows 33 and 35 - started second level modules
$ws.Range("J33").Value = $true
$ws.Range("J35").Value = $true

# Recalculate so dependent formulas (F33/F35/K33/K35/D8/D9/K39) refresh
$excel.Calculate()

# Update the view - scroll back up and select the last-edited note cell
[void]$ws.Range("H33").Select()
$excel.ActiveWindow.ScrollRow = 5
